$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new log entry (row 12) that was previously left blank.
# Write the Description cell first so the new shared string for the long
# description is registered before the shorter Subject string.
$ws.Range("D12").Value = "Finalized all grid and cell collection and also finalized the randomized `ndepth-first search with a stack implementation instead of recursive. `nHad issues with properly assigning the southern wall as my original `nwas very poorly done. Added a lot of comments across the project.`nLast step is a simple GUI to set the maze parameters and remake mazes."
$ws.Range("A12").Value = "Implemented algorithm and added comments"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = "5/23/2024"

# Copy the formatting from the row above (row 11) onto the new row so the
# newly populated cells look consistent with the rest of the log. D12 picks
# up the same literal style used by D10/D11.
$ws.Range("A11:D11").Copy() | Out-Null
$ws.Range("A12:D12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Rows.Item(12).RowHeight = 65.25

# Update the view state to match where the user left the selection/scroll.
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("G12").Select()
